# Regenerate save_data "K" column (column G) values with the recalculated
# "s_vals" (Strike# replaced by K), writing the new integer values in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$newValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    13 = 1
    14 = 3
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 2
    20 = 0
    21 = 2
    22 = 3
    23 = 0
    24 = 3
    25 = 1
    26 = 4
    27 = 1
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 2
    33 = 2
    34 = 0
    35 = 0
    36 = 1
    37 = 2
    38 = 2
    39 = 0
    40 = 1
    41 = 0
    42 = 2
    43 = 2
    44 = 0
    45 = 3
    46 = 1
    47 = 1
    48 = 0
    49 = 1
    50 = 1
    51 = 1
    52 = 2
    53 = 0
    54 = 1
    55 = 1
    56 = 6
    57 = 0
    58 = 2
    59 = 0
    60 = 1
    61 = 0
    63 = 1
    65 = 3
    66 = 2
    67 = 3
    68 = 1
    69 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
